$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E151").Value = "The release cycle is extremely long. It takes 6months or more to release software into production. The team needs to work on reducing the cycle time"
$ws.Range("E151").RowHeight = 43.2

$ws.Range("E152").Value = "Deployment is manual and very complex. It is generally a ceremony which requires high coordination between several stakeholders such as DBAs, Ops, Devs etc. This needs to be optimized and simplified so that the team is responsible for the deployment. Manual interventions should be minimized."
$ws.Range("E152").RowHeight = 72

$ws.Range("E153").Value = "The pace and delivery of work is inconsistent. The team is unable to give a definitive answer on when features will be delivered in production or how long it will take for the team to avail a certain feature to the customers upon completion of the development tasks."
$ws.Range("E153").RowHeight = 57.6

$ws.Range("E154").Value = "Each environment has it's own artefacts with high levels of configuration required to make them work in those environments. The team needs to start moving towards a configuration-as-code approach where the artefacts can simply be promoted from one environment to the other painlessly."
$ws.Range("E154").RowHeight = 72

$ws.Range("E155").Value = "The team is able to plan releases to a regular cycle and can track and manage features within regular time frames. Delivery is more predictable and captured on the team walls."
$ws.Range("E155").RowHeight = 43.2

$ws.Range("E156").Value = "The release dates are seldom moved out. The team is able to rather reduce the release features rather than the move out the release dates."
$ws.Range("E156").RowHeight = 28.8

$ws.Range("E157").Value = "There is a well defined process for handling the release of the artefacts to operations and support teams in places. This is documented and well understood by all the team members. The team at this level is investigating means of embedding operations/support in the team structure."
$ws.Range("E157").RowHeight = 72

$ws.Range("E158").Value = "Build artefacts are well defined and produced preferably by the build/continuous integration tool. These artefacts are readily available to the team as well as the operations and support teams. The team is able to tell which artefact was deployed when."
$ws.Range("E158").RowHeight = 57.6

$ws.Range("E159").Value = "The releases and deployments are automated using pipelines. These pipelines are visible on dashboards in the team spaces. There may however still be some manual dependencies such as environment configurations but the artefacts are packaged in a ready to deploy state."
$ws.Range("E159").RowHeight = 57.6

$ws.Range("E160").Value = "Teams have a high level of confidence for releases. Releases are no longer considered as a ceremonial activity but as a regular cadence."
$ws.Range("E160").RowHeight = 28.8

$ws.Range("E161").Value = "Release notes are automatically generated as part of the release/build scripts."
$ws.Range("E161").RowHeight = 28.8

$ws.Range("E162").Value = "The team does not have a key man dependency for releases. The release process is automated such that any individual in the team can perform a release. Releases can be audited and tracked."
$ws.Range("E162").RowHeight = 43.2

$ws.Range("E163").Value = "The operations teams collaborate very closely with the delivery teams and have backlog item tasks during development of the features to ensure that on completing development the applications are also ready to be deployed and released."
$ws.Range("E163").RowHeight = 57.6

$ws.Range("E164").Value = "The team constantly tracks and monitors the release metrics and trends on the story wall and all stakeholders are able to interact with these metrics on the dashboards."
$ws.Range("E164").RowHeight = 43.2

$ws.Range("E165").Value = "Releases can be deployed at anytime due to complete automation of the release tasks. Teams have release scripts and these are integrated in the build pipelines and dashboards."
$ws.Range("E165").RowHeight = 43.2

$ws.Range("E166").Value = "The team controls the release and does not depend on a third party team to compile and release their software. They can release software anytime to any of the environments."
$ws.Range("E166").RowHeight = 43.2

$ws.Range("E167").Value = "The product owner can decide when the release goes into production. In certain instances, the product owner can click on the deploy to production button and trigger a release to production."
$ws.Range("E167").RowHeight = 43.2

$ws.Range("E168").Value = "The operations team is part of the product/program team as opposed to being a hand-off team. They work hand-in-hand with the development teams to ensure the readiness of the application for production releases."
$ws.Range("E168").RowHeight = 43.2

$ws.Range("E169").Value = "Releases happen frequently and are routine. There are no application roll backs but rather roll forwards. Redeployment is favoured over system roll backs."

$ws.Range("E170").Value = "Team is able to release toggled off features using methodologies like Feature Toggling or Branching by Abstraction. Team is able to perform Blue-Green deployments with minimal downtimes to the application from a customer perspective."
$ws.Range("E170").RowHeight = 57.6

$ws.Range("E171").Select()
$excel.ActiveWindow.ScrollRow = 168
